$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.731.51"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "2.612.85"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'601.58"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").Value = "'154.49"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +1.61%  "

$ws.Range("D9").Value = "2.613.00"
$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("E10").Value = "  +7.37%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("E13").Value = "  -1.16%  "

$ws.Range("D14").Value = "'28.05"
$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("D16").Value = "3.087.55"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "67.614.82"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").Value = "2.611.20"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").Value = "'11.25"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").Value = "'366.01"
$ws.Range("E20").Value = "  +3.35%  "

$ws.Range("D21").Value = "'7.58"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("D23").Value = "'2.09"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'70.01"
$ws.Range("E25").Value = "  -0.64%  "

$ws.Range("D26").Value = "'10.05"
$ws.Range("E26").Value = "  -4.30%  "

$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("D29").Value = "'581.67"
$ws.Range("E29").Value = "  -2.71%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("D32").Value = "'7.93"
$ws.Range("E32").Value = "  -1.75%  "

$ws.Range("D33").Value = "'1.86"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("E36").Value = "  -2.54%  "

$ws.Range("D37").Value = "'4.96"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").Value = "'155.54"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").Value = "'5.39"
$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("D43").Value = "'2.65"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").Value = "'41.13"
$ws.Range("E44").Value = "  -0.68%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'16.43"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("D47").Value = "'156.37"
$ws.Range("E47").Value = "  -0.27%  "

$ws.Range("E48").Value = "  -7.74%  "

$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Value = "'20.97"
$ws.Range("E50").Value = "  -1.68%  "

$ws.Range("D51").Value = "'0.623"
$ws.Range("E51").Value = "  +0.34%  "
